$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.590.73'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.968.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.33%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.11%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.976.67'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.00%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.01%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -7.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.368'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.490.05'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.124'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.638.14'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.72'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.972.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.06%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -5.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '381.61'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.29%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.67'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.64%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.16'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.472'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.093.39'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.38%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0935'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -8.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.34'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.69%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.44'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '158.84'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.66'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.95'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.17%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.51%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.93'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.418.21'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -10.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.33'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '37.07'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.665'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0594'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.06%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0246'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.96'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -9.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0956'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.86%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.79'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.89%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '268.09'
